$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 25.5
$ws.Range("I8").Value = 25.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 76.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 62.5
$ws.Range("N8").ClearContents()

$ws.Range("H132").Value = 1806.7903
$ws.Range("I132").Value = 1471.8518
$ws.Range("K132").Value = 4415.555399999999
$ws.Range("M132").Value = -1885.555399999999

$ws.Range("H134").Value = 107673.336
$ws.Range("J134").Value = 107673.336
$ws.Range("L134").Value = 107673.336
$ws.Range("N134").Value = -117813.336


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 100000
$ws.Range("I10").Value = 100000
$ws.Range("K10").Value = 100000
$ws.Range("M10").Value = -99830

$ws.Range("H63").Value = 5117.7144
$ws.Range("I63").Value = 3838.7778
$ws.Range("J63").Value = 7419.8
$ws.Range("K63").Value = 3838.7778
$ws.Range("L63").Value = 7419.8
$ws.Range("M63").Value = -3152.7778
$ws.Range("N63").Value = -8791.799999999999

$ws.Range("H66").Value = 5117.7144
$ws.Range("I66").Value = 3838.7778
$ws.Range("J66").Value = 7419.8
$ws.Range("K66").Value = 19193.889
$ws.Range("L66").Value = 37099
$ws.Range("M66").Value = -15761.889
$ws.Range("N66").Value = -43963

$ws.Range("H92").Value = 160000
$ws.Range("J92").Value = 160000
$ws.Range("L92").Value = 160000
$ws.Range("N92").Value = -164992

$ws.Range("H139").Value = 89954.89
$ws.Range("J139").Value = 89954.89
$ws.Range("L139").Value = 89954.89
$ws.Range("N139").Value = -100234.89


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1270.6154
$ws.Range("I94").Value = 1059
$ws.Range("J94").Value = 1845
$ws.Range("K94").Value = 1059
$ws.Range("L94").Value = 1845
$ws.Range("M94").Value = -608
$ws.Range("N94").Value = -2747

$ws.Range("H105").Value = 3510.75
$ws.Range("I105").Value = 3510.75
$ws.Range("K105").Value = 3510.75
$ws.Range("M105").Value = -1763.75

$ws.Range("H135").Value = 53986.152
$ws.Range("J135").Value = 53986.152
$ws.Range("L135").Value = 53986.152
$ws.Range("N135").Value = -64126.152

$ws.Range("H137").Value = 49808.332
$ws.Range("J137").Value = 49808.332
$ws.Range("L137").Value = 49808.332
$ws.Range("N137").Value = -60008.332


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 74072.86
$ws.Range("I62").Value = 85951.664
$ws.Range("J62").Value = 2800
$ws.Range("K62").Value = 85951.664
$ws.Range("L62").Value = 2800
$ws.Range("M62").Value = -85327.664
$ws.Range("N62").Value = -4048

$ws.Range("H65").Value = 74072.86
$ws.Range("I65").Value = 85951.664
$ws.Range("J65").Value = 2800
$ws.Range("K65").Value = 429758.32
$ws.Range("L65").Value = 14000
$ws.Range("M65").Value = -426638.32
$ws.Range("N65").Value = -20240

$ws.Range("H110").Value = 25044.5
$ws.Range("J110").Value = 25044.5
$ws.Range("L110").Value = 25044.5
$ws.Range("N110").Value = -33224.5

$ws.Range("H135").Value = 85107.73
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 85107.73
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 85107.73
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -95247.73

$ws.Range("H138").Value = 52397
$ws.Range("J138").Value = 52397
$ws.Range("L138").Value = 52397
$ws.Range("N138").Value = -62677


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 238.35715
$ws.Range("J11").Value = 398.33334
$ws.Range("L11").Value = 1195.00002
$ws.Range("N11").Value = -1475.00002

$ws.Range("H134").Value = 4021.5454
$ws.Range("I134").Value = 1939.9048
$ws.Range("J134").Value = 7664.4165
$ws.Range("K134").Value = 5819.7144
$ws.Range("L134").Value = 22993.2495
$ws.Range("M134").Value = -749.7143999999998
$ws.Range("N134").Value = -33133.24950000001


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6429.9414
$ws.Range("I70").Value = 5300
$ws.Range("J70").Value = 6900.75
$ws.Range("K70").Value = 5300
$ws.Range("L70").Value = 6900.75
$ws.Range("M70").Value = -5030
$ws.Range("N70").Value = -7440.75

$ws.Range("H73").Value = 6429.9414
$ws.Range("I73").Value = 5300
$ws.Range("J73").Value = 6900.75
$ws.Range("K73").Value = 5300
$ws.Range("L73").Value = 6900.75
$ws.Range("M73").Value = -4364
$ws.Range("N73").Value = -8772.75

$ws.Range("H80").Value = 2937.7856
$ws.Range("I80").Value = 2875.4546
$ws.Range("J80").Value = 3166.3333
$ws.Range("K80").Value = 2875.4546
$ws.Range("L80").Value = 3166.3333
$ws.Range("M80").Value = -1877.4546
$ws.Range("N80").Value = -5162.3333

$ws.Range("H83").Value = 2937.7856
$ws.Range("I83").Value = 2875.4546
$ws.Range("J83").Value = 3166.3333
$ws.Range("K83").Value = 14377.273
$ws.Range("L83").Value = 15831.6665
$ws.Range("M83").Value = -9385.273000000001
$ws.Range("N83").Value = -25815.6665

$ws.Range("H132").Value = 2642.6667
$ws.Range("I132").Value = 2451.2222
$ws.Range("J132").Value = 2929.8333
$ws.Range("K132").Value = 7353.6666
$ws.Range("L132").Value = 8789.499899999999
$ws.Range("M132").Value = -4823.6666
$ws.Range("N132").Value = -13849.4999

$ws.Range("H135").Value = 48680.332
$ws.Range("J135").Value = 48680.332
$ws.Range("L135").Value = 48680.332
$ws.Range("N135").Value = -58820.332

$ws.Range("H138").Value = 48709.6
$ws.Range("J138").Value = 48709.6
$ws.Range("L138").Value = 48709.6
$ws.Range("N138").Value = -58989.6

$ws.Range("H140").Value = 49552.223
$ws.Range("J140").Value = 49552.223
$ws.Range("L140").Value = 49552.223
$ws.Range("N140").Value = -59912.223


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 23100
$ws.Range("J10").Value = 24650
$ws.Range("L10").Value = 24650
$ws.Range("N10").Value = -24930

$ws.Range("H134").Value = 75076.336
$ws.Range("J134").Value = 75076.336
$ws.Range("L134").Value = 75076.336
$ws.Range("N134").Value = -85216.336

$ws.Range("H136").Value = 3443.2744
$ws.Range("I136").Value = 3579.9333
$ws.Range("J136").Value = 2418.3333
$ws.Range("K136").Value = 10739.7999
$ws.Range("L136").Value = 7254.999899999999
$ws.Range("M136").Value = -8189.7999
$ws.Range("N136").Value = -12354.9999

$ws.Range("H137").Value = 72724.39999999999
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 72724.39999999999
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 72724.39999999999
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -82924.39999999999

$ws.Range("H139").Value = 37472.668
$ws.Range("J139").Value = 37472.668
$ws.Range("L139").Value = 37472.668
$ws.Range("N139").Value = -47752.668

$ws.Range("H141").Value = 45085
$ws.Range("J141").Value = 45085
$ws.Range("L141").Value = 45085
$ws.Range("N141").Value = -55445


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 59071.555
$ws.Range("J135").Value = 59071.555
$ws.Range("L135").Value = 59071.555
$ws.Range("N135").Value = -69211.55499999999

$ws.Range("H139").Value = 57625
$ws.Range("J139").Value = 57625
$ws.Range("L139").Value = 57625
$ws.Range("N139").Value = -67905

$ws.Range("H141").Value = 82289.5
$ws.Range("J141").Value = 82289.5
$ws.Range("L141").Value = 82289.5
$ws.Range("N141").Value = -92649.5

